$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3825.7144
$ws.Range("J19").Value = 4632.625
$ws.Range("L19").Value = 4632.625
$ws.Range("N19").Value = -4982.625
$ws.Range("H33").Value = 50926.6
$ws.Range("I33").Value = 68719
$ws.Range("K33").Value = 68719
$ws.Range("M33").Value = -68490
$ws.Range("H40").Value = 2230.3076
$ws.Range("I40").Value = 1824.75
$ws.Range("J40").Value = 2410.5557
$ws.Range("K40").Value = 1824.75
$ws.Range("L40").Value = 2410.5557
$ws.Range("M40").Value = -1649.75
$ws.Range("N40").Value = -2760.5557
$ws.Range("H86").Value = 5586370.5
$ws.Range("I86").Value = 2688.5334
$ws.Range("J86").Value = 9574715
$ws.Range("K86").Value = 2688.5334
$ws.Range("L86").Value = 9574715
$ws.Range("M86").Value = -1565.5334
$ws.Range("N86").Value = -9576961
$ws.Range("H87").Value = 128666.664
$ws.Range("J87").Value = 128666.664
$ws.Range("L87").Value = 128666.664
$ws.Range("N87").Value = -131162.664
$ws.Range("H89").Value = 5586370.5
$ws.Range("I89").Value = 2688.5334
$ws.Range("J89").Value = 9574715
$ws.Range("K89").Value = 13442.667
$ws.Range("L89").Value = 47873575
$ws.Range("M89").Value = -7826.666999999999
$ws.Range("N89").Value = -47884807
$ws.Range("H90").Value = 128666.664
$ws.Range("J90").Value = 128666.664
$ws.Range("L90").Value = 385999.992
$ws.Range("N90").Value = -398479.992
$ws.Range("H92").Value = 90691.05
$ws.Range("J92").Value = 199126.3
$ws.Range("L92").Value = 199126.3
$ws.Range("N92").Value = -201622.3
$ws.Range("H100").Value = 1441.2174
$ws.Range("I100").Value = 1223.7059
$ws.Range("J100").Value = 2057.5
$ws.Range("K100").Value = 1223.7059
$ws.Range("L100").Value = 2057.5
$ws.Range("M100").Value = -682.7058999999999
$ws.Range("N100").Value = -3139.5
$ws.Range("H132").Value = 2269.1528
$ws.Range("I132").Value = 2072.842
$ws.Range("J132").Value = 3926.889
$ws.Range("K132").Value = 6218.526
$ws.Range("L132").Value = 11780.667
$ws.Range("M132").Value = -3688.526
$ws.Range("N132").Value = -16840.667
$ws.Range("H138").Value = 2146.375
$ws.Range("I138").Value = 1837.591
$ws.Range("J138").Value = 2825.7
$ws.Range("K138").Value = 5512.772999999999
$ws.Range("L138").Value = 8477.099999999999
$ws.Range("M138").Value = -372.7729999999992
$ws.Range("N138").Value = -18757.1
$ws.Range("H141").Value = 975.4828
$ws.Range("I141").Value = 975.4828
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2926.4484
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2253.5516
$ws.Range("N141").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H32").Value = 3948
$ws.Range("I32").Value = 4567.56
$ws.Range("J32").Value = 1366.5
$ws.Range("K32").Value = 4567.56
$ws.Range("L32").Value = 1366.5
$ws.Range("M32").Value = -4280.56
$ws.Range("N32").Value = -1940.5
$ws.Range("H132").Value = 1324.9
$ws.Range("I132").Value = 1245.4615
$ws.Range("J132").Value = 1472.4286
$ws.Range("K132").Value = 3736.3845
$ws.Range("L132").Value = 4417.2858
$ws.Range("M132").Value = -1206.3845
$ws.Range("N132").Value = -9477.2858

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 644.16
$ws.Range("I22").Value = 586.95
$ws.Range("J22").Value = 873
$ws.Range("K22").Value = 586.95
$ws.Range("L22").Value = 873
$ws.Range("M22").Value = -236.95
$ws.Range("N22").Value = -1573
$ws.Range("H31").Value = 6002.6597
$ws.Range("I31").Value = 14502.454
$ws.Range("K31").Value = 14502.454
$ws.Range("M31").Value = -14207.454
$ws.Range("H34").Value = 6002.6597
$ws.Range("I34").Value = 14502.454
$ws.Range("K34").Value = 14502.454
$ws.Range("M34").Value = -14300.454
$ws.Range("H58").Value = 1420.9302
$ws.Range("I58").Value = 1347.0322
$ws.Range("J58").Value = 1611.8334
$ws.Range("K58").Value = 1347.0322
$ws.Range("L58").Value = 1611.8334
$ws.Range("M58").Value = -1144.0322
$ws.Range("N58").Value = -2017.8334
$ws.Range("H59").Value = 49124.5
$ws.Range("I59").Value = 50874
$ws.Range("J59").Value = 47375
$ws.Range("K59").Value = 50874
$ws.Range("L59").Value = 47375
$ws.Range("M59").Value = -49729
$ws.Range("N59").Value = -49665
$ws.Range("H60").Value = 24916.666
$ws.Range("I60").Value = 24500
$ws.Range("K60").Value = 24500
$ws.Range("M60").Value = -23989
$ws.Range("H130").Value = 75000
$ws.Range("J130").Value = 75000
$ws.Range("L130").Value = 75000
$ws.Range("N130").Value = -85040
$ws.Range("H132").Value = 3352.2334
$ws.Range("I132").Value = 2923.1365
$ws.Range("K132").Value = 8769.4095
$ws.Range("M132").Value = -6239.4095
$ws.Range("H134").Value = 1662.919
$ws.Range("I134").Value = 1534.8438
$ws.Range("K134").Value = 4604.5314
$ws.Range("M134").Value = -2069.5314
$ws.Range("H136").Value = 1420.9302
$ws.Range("I136").Value = 1347.0322
$ws.Range("J136").Value = 1611.8334
$ws.Range("K136").Value = 4041.0966
$ws.Range("L136").Value = 4835.5002
$ws.Range("M136").Value = -1491.0966
$ws.Range("N136").Value = -9935.5002

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1014.3684
$ws.Range("I5").Value = 622.9
$ws.Range("J5").Value = 1449.3334
$ws.Range("K5").Value = 1868.7
$ws.Range("L5").Value = 4348.0002
$ws.Range("M5").Value = -1756.7
$ws.Range("N5").Value = -4572.0002
$ws.Range("H26").Value = 593.5
$ws.Range("I26").Value = 594
$ws.Range("J26").Value = 593
$ws.Range("K26").Value = 1782
$ws.Range("L26").Value = 1779
$ws.Range("M26").Value = -1494
$ws.Range("N26").Value = -2355
$ws.Range("H37").Value = 38750
$ws.Range("J37").Value = 38750
$ws.Range("L37").Value = 116250
$ws.Range("N37").Value = -116474
$ws.Range("H135").Value = 1014.3684
$ws.Range("I135").Value = 622.9
$ws.Range("J135").Value = 1449.3334
$ws.Range("K135").Value = 5606.099999999999
$ws.Range("L135").Value = 13044.0006
$ws.Range("M135").Value = -3071.099999999999
$ws.Range("N135").Value = -18114.0006

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3952.5
$ws.Range("I43").Value = 2311.5386
$ws.Range("J43").Value = 7000
$ws.Range("K43").Value = 2311.5386
$ws.Range("L43").Value = 7000
$ws.Range("M43").Value = -2160.5386
$ws.Range("N43").Value = -7302
$ws.Range("H46").Value = 7167.364
$ws.Range("I46").Value = 2884.1
$ws.Range("K46").Value = 2884.1
$ws.Range("M46").Value = -2728.1
$ws.Range("H57").Value = 23189
$ws.Range("I57").Value = 10236.25
$ws.Range("K57").Value = 10236.25
$ws.Range("M57").Value = -9416.25
$ws.Range("H70").Value = 9778.714
$ws.Range("J70").Value = 7876.25
$ws.Range("L70").Value = 7876.25
$ws.Range("N70").Value = -8416.25
$ws.Range("H73").Value = 9778.714
$ws.Range("J73").Value = 7876.25
$ws.Range("L73").Value = 7876.25
$ws.Range("N73").Value = -9748.25
$ws.Range("H80").Value = 3997.682
$ws.Range("I80").Value = 3888.9285
$ws.Range("J80").Value = 4188
$ws.Range("K80").Value = 3888.9285
$ws.Range("L80").Value = 4188
$ws.Range("M80").Value = -2890.9285
$ws.Range("N80").Value = -6184
$ws.Range("H83").Value = 3997.682
$ws.Range("I83").Value = 3888.9285
$ws.Range("J83").Value = 4188
$ws.Range("K83").Value = 19444.6425
$ws.Range("L83").Value = 20940
$ws.Range("M83").Value = -14452.6425
$ws.Range("N83").Value = -30924
$ws.Range("H102").Value = 2611.8635
$ws.Range("I102").Value = 1399.6666
$ws.Range("J102").Value = 4066.5
$ws.Range("K102").Value = 1399.6666
$ws.Range("L102").Value = 4066.5
$ws.Range("M102").Value = 222.3334
$ws.Range("N102").Value = -7310.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2608.6775
$ws.Range("I46").Value = 1529.5
$ws.Range("J46").Value = 3759.8
$ws.Range("K46").Value = 1529.5
$ws.Range("L46").Value = 3759.8
$ws.Range("M46").Value = -1341.5
$ws.Range("N46").Value = -4135.8
$ws.Range("H132").Value = 4085.4167
$ws.Range("I132").Value = 2531.3333
$ws.Range("J132").Value = 6675.5557
$ws.Range("K132").Value = 7593.999899999999
$ws.Range("L132").Value = 20026.6671
$ws.Range("M132").Value = -5063.999899999999
$ws.Range("N132").Value = -25086.6671

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 46998
$ws.Range("J93").Value = 46998
$ws.Range("L93").Value = 46998
$ws.Range("N93").Value = -51990
$ws.Range("H122").Value = 1336.7567
$ws.Range("I122").Value = 1162.5333
$ws.Range("J122").Value = 2083.4285
$ws.Range("K122").Value = 3487.5999
$ws.Range("L122").Value = 6250.2855
$ws.Range("M122").Value = -1037.5999
$ws.Range("N122").Value = -11150.2855
$ws.Range("H136").Value = 2977.652
$ws.Range("I136").Value = 1192.3334
$ws.Range("J136").Value = 6325.125
$ws.Range("K136").Value = 3577.0002
$ws.Range("L136").Value = 18975.375
$ws.Range("M136").Value = -1027.0002
$ws.Range("N136").Value = -24075.375
